$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("G4").Value = 1.75
$ws.Range("H4").Value = 3.4
$ws.Range("I4").Value = 5
$ws.Range("N4").Value = 8
$ws.Range("Z4").Value = 13
$ws.Range("AH4").Value = 11
$ws.Range("AI4").Value = 23
$ws.Range("AJ4").Value = 17
$ws.Range("G5").Value = 2.6
$ws.Range("H5").Value = 2.88
$ws.Range("G8").Value = 4.2
$ws.Range("H8").Value = 3.4
$ws.Range("I8").Value = 1.73
$ws.Range("K8").Value = 2.2
$ws.Range("L8").Value = 2.4
$ws.Range("AA8").Value = 34
$ws.Range("AK8").Value = 15
$ws.Range("AL8").Value = 15
$ws.Range("AN8").Value = 6
$ws.Range("AX8").Value = 9.5
$ws.Range("H11").Value = 4.33
$ws.Range("I11").Value = 4.75
$ws.Range("J11").Value = 2.1
$ws.Range("L11").Value = 4.75
$ws.Range("Q11").Value = 1.53
$ws.Range("R11").Value = 2.4
$ws.Range("S11").Value = 1.29
$ws.Range("T11").Value = 3.5
$ws.Range("W11").Value = 9.5
$ws.Range("X11").Value = 9.5
$ws.Range("Z11").Value = 13
$ws.Range("AB11").Value = 21
$ws.Range("AI11").Value = 29
$ws.Range("AL11").Value = 34
$ws.Range("AM11").Value = 34
$ws.Range("AO11").Value = 8
$ws.Range("AT11").Value = 3.5
$ws.Range("AW11").Value = 7
$ws.Range("AX11").Value = 23
$ws.Range("AY11").Value = 26
$ws.Range("BA11").Value = 81
$ws.Range("BB11").Value = 151
$ws.Range("G12").Value = 2.15
$ws.Range("I12").Value = 3.2
$ws.Range("J12").Value = 2.75
$ws.Range("L12").Value = 3.5
$ws.Range("M12").Value = 1.03
$ws.Range("N12").Value = 15
$ws.Range("S12").Value = 1.3
$ws.Range("T12").Value = 3.4
$ws.Range("U12").Value = 1.53
$ws.Range("V12").Value = 2.38
$ws.Range("W12").Value = 10
$ws.Range("X12").Value = 12
$ws.Range("Z12").Value = 21
$ws.Range("AC12").Value = 15
$ws.Range("AF12").Value = 34
$ws.Range("AI12").Value = 19
$ws.Range("AJ12").Value = 12
$ws.Range("AK12").Value = 34
$ws.Range("AL12").Value = 23
$ws.Range("AM12").Value = 26
$ws.Range("AN12").Value = 4.5
$ws.Range("AT12").Value = 3.4
$ws.Range("AU12").Value = 7
$ws.Range("AX12").Value = 17
$ws.Range("AY12").Value = 21
$ws.Range("BA12").Value = 51
$ws.Range("U13").Value = 1.5
$ws.Range("V13").Value = 2.5
$ws.Range("AC13").Value = 21
$ws.Range("AG13").Value = 101
$ws.Range("AJ13").Value = 8.5
$ws.Range("AP13").Value = 23
$ws.Range("AW13").Value = 4
$ws.Range("G14").Value = 4.5
$ws.Range("I14").Value = 1.57
$ws.Range("J14").Value = 4.5
$ws.Range("K14").Value = 2.75
$ws.Range("L14").Value = 2.05
$ws.Range("U14").Value = 1.44
$ws.Range("V14").Value = 2.63
$ws.Range("W14").Value = 23
$ws.Range("Z14").Value = 51
$ws.Range("AA14").Value = 29
$ws.Range("AB14").Value = 29
$ws.Range("AE14").Value = 13
$ws.Range("AF14").Value = 34
$ws.Range("AI14").Value = 11
$ws.Range("AL14").Value = 11
$ws.Range("AQ14").Value = 67
$ws.Range("AR14").Value = 67
$ws.Range("AS14").Value = 101
$ws.Range("AV14").Value = 41
$ws.Range("AX14").Value = 8
$ws.Range("BC14").Value = 251
$ws.Range("U22").Value = 1.57
$ws.Range("V23").Value = 1.73
$ws.Range("U24").Value = 1.91
$ws.Range("V24").Value = 1.8
$ws.Range("G28").Value = 1.85
$ws.Range("H28").Value = 3.4
$ws.Range("I28").Value = 4.33
$ws.Range("J28").Value = 2.5
$ws.Range("L28").Value = 4.5
$ws.Range("Q28").Value = 1.98
$ws.Range("R28").Value = 1.88
$ws.Range("Z28").Value = 15
$ws.Range("BA28").Value = 101
$ws.Range("G30").Value = 2
$ws.Range("I30").Value = 3.8
$ws.Range("L30").Value = 4.33
$ws.Range("Q30").Value = 2.05
$ws.Range("R30").Value = 1.75
$ws.Range("S30").Value = 1.44
$ws.Range("T30").Value = 2.63
$ws.Range("U30").Value = 1.83
$ws.Range("V30").Value = 1.83
$ws.Range("W30").Value = 7
$ws.Range("X30").Value = 9.5
$ws.Range("Z30").Value = 17
$ws.Range("AL30").Value = 34
$ws.Range("AP30").Value = 23
$ws.Range("AT30").Value = 2.63
$ws.Range("AX30").Value = 21
$ws.Range("BA30").Value = 101
$ws.Range("M38").Value = 1.08
$ws.Range("O38").Value = 1.44
$ws.Range("P38").Value = 2.63
$ws.Range("M39").Value = 1.05
$ws.Range("O39").Value = 1.29
